# Add team record columns (Wins, Losses, Ties) to the roster sheet.
# New columns AD, AE, AF are appended after the existing data (A:AC),
# with header labels in row 1 and a constant 81-81-0 record for every
# player row (2-55), matching the website_scraper style addition.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colWins   = 30  # AD
$colLosses = 31  # AE
$colTies   = 32  # AF

$lastRow = 55

# Header row (row 1) - add labels and copy the existing header formatting
# (bold font, border, centered alignment) from A1 onto the new header cells.
$ws.Cells.Item(1, $colWins).Value   = "Wins"
$ws.Cells.Item(1, $colLosses).Value = "Losses"
$ws.Cells.Item(1, $colTies).Value   = "Ties"

$ws.Cells.Item(1, 1).Copy()
$headerDest = $ws.Range($ws.Cells.Item(1, $colWins), $ws.Cells.Item(1, $colTies))
$headerDest.PasteSpecial(-4122)

# Data rows - constant team record (81 wins, 81 losses, 0 ties) for every row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $colWins).Value   = 81
    $ws.Cells.Item($r, $colLosses).Value = 81
    $ws.Cells.Item($r, $colTies).Value   = 0
}
